$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304, shifting existing rows 304-408 down to 305-409
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new data record
$ws.Range("A304").Value = 5
$ws.Range("B304").Value = "Macroferia Regional de Talca"
$ws.Range("C304").Value = "Maule"
$ws.Range("D304").Value = 44559
$ws.Range("E304").Value = 7
$ws.Range("F304").Value = 100112002
$ws.Range("G304").Value = "Pimiento"
$ws.Range("H304").Value = "Cuatro cascos verde"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 300
$ws.Range("K304").Value = 11000
$ws.Range("L304").Value = 11000
$ws.Range("M304").Value = 11000
$ws.Range("N304").Value = "$/caja 15 kilos"
$ws.Range("O304").Value = "Región del Maule"
$ws.Range("P304").Value = 733
$ws.Range("Q304").Value = 15
$ws.Range("R304").Value = "Hortaliza"
